$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$ws.Range("G2:I2").NumberFormat = "@"
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "456"
$ws.Range("I2").Value = "5"
